$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Use Case")

# Move the value currently in C2 down to C7 (cut/paste),
# then put the new name in C2.
$ws.Range("C2").Cut($ws.Range("C7"))
$ws.Range("C2").Value = "Эвертов Владимир Васильевич"

$ws.Range("C7").Select()
